$d = $word.ActiveDocument

$replacements = @(
    @{old="955×2=1910"; new="230×3=690"},
    @{old="977×4=3908"; new="121×3=363"},
    @{old="820×7=5740"; new="452×4=1808"},
    @{old="553×9=4977"; new="130×5=650"},
    @{old="827×4=3308"; new="349×7=2443"},
    @{old="482×7=3374"; new="219×3=657"},
    @{old="470×4=1880"; new="507×8=4056"},
    @{old="375×7=2625"; new="811×7=5677"},
    @{old="591×9=5319"; new="745×7=5215"},
    @{old="430×3=1290"; new="461×3=1383"},
    @{old="842×2=1684"; new="142×7=994"},
    @{old="975×7=6825"; new="299×3=897"},
    @{old="662×3=1986"; new="492×9=4428"},
    @{old="579×5=2895"; new="323×3=969"},
    @{old="438×4=1752"; new="501×9=4509"},
    @{old="914×3=2742"; new="878×2=1756"},
    @{old="353×7=2471"; new="326×8=2608"},
    @{old="456×6=2736"; new="384×4=1536"},
    @{old="626×6=3756"; new="985×8=7880"},
    @{old="250×8=2000"; new="659×7=4613"},
    @{old="927×4=3708"; new="459×7=3213"},
    @{old="266×3=798";  new="380×5=1900"},
    @{old="492×3=1476"; new="846×9=7614"},
    @{old="877×5=4385"; new="696×6=4176"},
    @{old="457×7=3199"; new="423×4=1692"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
